# Update forecast values (column C) on Sheet1 to reflect refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 313.2
    3  = 315.7
    4  = 324.6
    5  = 325.5
    6  = 314.8
    7  = 313.5
    8  = 324.1
    9  = 325.7
    10 = 327.4
    11 = 326.3
    12 = 329.1
    13 = 327.6
    14 = 325.3
    15 = 317.1
    16 = 316.4
    17 = 300.6
    18 = 301.3
    19 = 305
    20 = 299.1
    21 = 294.5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
